$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "Pasien A"
$ws.Cells.Item(2, 2).Value = "L"
$ws.Cells.Item(2, 3).Value = 9
$ws.Cells.Item(2, 4).Value = "A"
$ws.Cells.Item(2, 5).Value = "Demam Gastritis"
$ws.Cells.Item(2, 6).Value = "B70 K29"

$ws.Cells.Item(3, 1).Value = "Pasien B"
$ws.Cells.Item(3, 2).Value = "P"
$ws.Cells.Item(3, 3).Value = 17
$ws.Cells.Item(3, 4).Value = "B"
$ws.Cells.Item(3, 5).Value = "Sakit_Kepala"
$ws.Cells.Item(3, 6).Value = "G30"

$ws.Cells.Item(4, 1).Value = "Pasien C"
$ws.Cells.Item(4, 2).Value = "L"
$ws.Cells.Item(4, 3).Value = 51
$ws.Cells.Item(4, 4).Value = "X"
$ws.Cells.Item(4, 5).Value = "Diare Demam"
$ws.Cells.Item(4, 6).Value = "K10 B70"

$ws.Cells.Item(5, 1).Value = "Pasien D"
$ws.Cells.Item(5, 2).Value = "P"
$ws.Cells.Item(5, 3).Value = 36
$ws.Cells.Item(5, 4).Value = "Y"
$ws.Cells.Item(5, 5).Value = "Sakit_Pinggang"
$ws.Cells.Item(5, 6).Value = "G20"

$ws.Cells.Item(6, 1).Value = "Pasien E"
$ws.Cells.Item(6, 2).Value = "L"
$ws.Cells.Item(6, 3).Value = 25
$ws.Cells.Item(6, 4).Value = "Z"
$ws.Cells.Item(6, 5).Value = "Flu"
$ws.Cells.Item(6, 6).Value = "B50"

$ws.Cells.Item(7, 1).Value = "Pasien F"
$ws.Cells.Item(7, 2).Value = "L"
$ws.Cells.Item(7, 3).Value = 70
$ws.Cells.Item(7, 4).Value = "H"
$ws.Cells.Item(7, 5).Value = "Anemia"
$ws.Cells.Item(7, 6).Value = "D50"

$ws.Cells.Item(8, 1).Value = "Pasien G"
$ws.Cells.Item(8, 2).Value = "L"
$ws.Cells.Item(8, 3).Value = 18
$ws.Cells.Item(8, 4).Value = "A"
$ws.Cells.Item(8, 5).Value = "Demam"
$ws.Cells.Item(8, 6).Value = "B70"

$ws.Cells.Item(9, 1).Value = "Pasien H"
$ws.Cells.Item(9, 2).Value = "P"
$ws.Cells.Item(9, 3).Value = 41
$ws.Cells.Item(9, 4).Value = "B"
$ws.Cells.Item(9, 5).Value = "Sakit_Kepala"
$ws.Cells.Item(9, 6).Value = "G30"

$ws.Cells.Item(10, 1).Value = "Pasien I"
$ws.Cells.Item(10, 2).Value = "P"
$ws.Cells.Item(10, 3).Value = 23
$ws.Cells.Item(10, 4).Value = "X"
$ws.Cells.Item(10, 5).Value = "Diare Demam"
$ws.Cells.Item(10, 6).Value = "K10 B70"

$ws.Cells.Item(11, 1).Value = "Pasien J"
$ws.Cells.Item(11, 2).Value = "P"
$ws.Cells.Item(11, 3).Value = 55
$ws.Cells.Item(11, 4).Value = "Y"
$ws.Cells.Item(11, 5).Value = "Sakit_Pinggang"
$ws.Cells.Item(11, 6).Value = "G20"

$ws.Cells.Item(12, 1).Value = "Pasien K"
$ws.Cells.Item(12, 2).Value = "L"
$ws.Cells.Item(12, 3).Value = 2
$ws.Cells.Item(12, 4).Value = "Z"
$ws.Cells.Item(12, 5).Value = "Flu"
$ws.Cells.Item(12, 6).Value = "B50"

$ws.Cells.Item(13, 1).Value = "Pasien L"
$ws.Cells.Item(13, 2).Value = "P"
$ws.Cells.Item(13, 3).Value = 38
$ws.Cells.Item(13, 4).Value = "H"
$ws.Cells.Item(13, 5).Value = "Anemia"
$ws.Cells.Item(13, 6).Value = "D50"

$ws.Cells.Item(14, 1).Value = "Pasien M"
$ws.Cells.Item(14, 2).Value = "L"
$ws.Cells.Item(14, 3).Value = 34
$ws.Cells.Item(14, 4).Value = "I"
$ws.Cells.Item(14, 5).Value = "Demam"
$ws.Cells.Item(14, 6).Value = "B70"

$ws.Cells.Item(15, 1).Value = "Pasien N"
$ws.Cells.Item(15, 2).Value = "P"
$ws.Cells.Item(15, 3).Value = 56
$ws.Cells.Item(15, 4).Value = "J"
$ws.Cells.Item(15, 5).Value = "Sakit_Kepala"
$ws.Cells.Item(15, 6).Value = "G30"

$ws.Cells.Item(16, 1).Value = "Pasien O"
$ws.Cells.Item(16, 2).Value = "L"
$ws.Cells.Item(16, 3).Value = 71
$ws.Cells.Item(16, 4).Value = "K"
$ws.Cells.Item(16, 5).Value = "Demam DBD"
$ws.Cells.Item(16, 6).Value = "B70 C10"

$ws.Cells.Item(17, 1).Value = "Pasien P"
$ws.Cells.Item(17, 2).Value = "P"
$ws.Cells.Item(17, 3).Value = 23
$ws.Cells.Item(17, 4).Value = "L"
$ws.Cells.Item(17, 5).Value = "DBD"
$ws.Cells.Item(17, 6).Value = "C10"

$ws.Cells.Item(18, 1).Value = "Pasien Q"
$ws.Cells.Item(18, 2).Value = "L"
$ws.Cells.Item(18, 3).Value = 11
$ws.Cells.Item(18, 4).Value = "M"
$ws.Cells.Item(18, 5).Value = "Sakit_Pinggang"
$ws.Cells.Item(18, 6).Value = "G20"

$ws.Cells.Item(19, 1).Value = "Pasien R"
$ws.Cells.Item(19, 2).Value = "L"
$ws.Cells.Item(19, 3).Value = 0.5
$ws.Cells.Item(19, 4).Value = "N"
$ws.Cells.Item(19, 5).Value = "Flu"
$ws.Cells.Item(19, 6).Value = "B50"

$ws.Cells.Item(20, 1).Value = "Pasien S"
$ws.Cells.Item(20, 2).Value = "L"
$ws.Cells.Item(20, 3).Value = 28
$ws.Cells.Item(20, 4).Value = "O"
$ws.Cells.Item(20, 5).Value = "Anemia Gastritis"
$ws.Cells.Item(20, 6).Value = "D50 K29"

$ws.Cells.Item(21, 1).Value = "Pasien T"
$ws.Cells.Item(21, 2).Value = "P"
$ws.Cells.Item(21, 3).Value = 39
$ws.Cells.Item(21, 4).Value = "P"
$ws.Cells.Item(21, 5).Value = "Sakit_Gigi Demam"
$ws.Cells.Item(21, 6).Value = "G40  B70"

$ws.Cells.Item(22, 1).Value = "Pasien U"
$ws.Cells.Item(22, 2).Value = "P"
$ws.Cells.Item(22, 3).Value = 45
$ws.Cells.Item(22, 4).Value = "Q"
$ws.Cells.Item(22, 5).Value = "Gastritis"
$ws.Cells.Item(22, 6).Value = "K29"

$ws.Cells.Item(23, 1).Value = "Pasien V"
$ws.Cells.Item(23, 2).Value = "L"
$ws.Cells.Item(23, 3).Value = 51
$ws.Cells.Item(23, 4).Value = "R"
$ws.Cells.Item(23, 5).Value = "Demam"
$ws.Cells.Item(23, 6).Value = "B70"

$ws.Cells.Item(24, 1).Value = "Pasien W"
$ws.Cells.Item(24, 2).Value = "L"
$ws.Cells.Item(24, 3).Value = 12
$ws.Cells.Item(24, 4).Value = "S"
$ws.Cells.Item(24, 5).Value = "Sakit_Kepala"
$ws.Cells.Item(24, 6).Value = "G30"

$ws.Cells.Item(25, 1).Value = "Pasien X"
$ws.Cells.Item(25, 2).Value = "L"
$ws.Cells.Item(25, 3).Value = 18
$ws.Cells.Item(25, 4).Value = "T"
$ws.Cells.Item(25, 5).Value = "Anemia Gastritis"
$ws.Cells.Item(25, 6).Value = "D50 K29"

$ws.Cells.Item(26, 1).Value = "Pasien Y"
$ws.Cells.Item(26, 2).Value = "P"
$ws.Cells.Item(26, 3).Value = 48
$ws.Cells.Item(26, 4).Value = "U"
$ws.Cells.Item(26, 5).Value = "Demam Sakit_Kepala"
$ws.Cells.Item(26, 6).Value = "B70 G30"

$ws.Cells.Item(27, 1).Value = "Pasien Z"
$ws.Cells.Item(27, 2).Value = "L"
$ws.Cells.Item(27, 3).Value = 21
$ws.Cells.Item(27, 4).Value = "V"
$ws.Cells.Item(27, 5).Value = "Sakit_Kepala"
$ws.Cells.Item(27, 6).Value = "G30"

$ws.Cells.Item(28, 1).Value = "Pasien AA"
$ws.Cells.Item(28, 2).Value = "P"
$ws.Cells.Item(28, 3).Value = 75
$ws.Cells.Item(28, 4).Value = "W"
$ws.Cells.Item(28, 5).Value = "Diare Demam"
$ws.Cells.Item(28, 6).Value = "K10 B70"

$ws.Cells.Item(29, 1).Value = "Pasien AB"
$ws.Cells.Item(29, 2).Value = "L"
$ws.Cells.Item(29, 3).Value = 63
$ws.Cells.Item(29, 4).Value = "X"
$ws.Cells.Item(29, 5).Value = "Sakit_Pinggang"
$ws.Cells.Item(29, 6).Value = "G20"

$ws.Cells.Item(30, 1).Value = "Pasien AC"
$ws.Cells.Item(30, 2).Value = "P"
$ws.Cells.Item(30, 3).Value = 58
$ws.Cells.Item(30, 4).Value = "Y"
$ws.Cells.Item(30, 5).Value = "Sakit_Gigi"
$ws.Cells.Item(30, 6).Value = "G40"

$ws.Columns.Item(5).ColumnWidth = 17.5
[void]$ws.Range("F24").Select()
